$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "El plan infinito"
$ws.Range("B8").Value = "Isabel Allende"
$ws.Range("C8").Value = ""
